# Refresh BCH yearly financial statement figures (income statement, balance
# sheet, and cash flow statement) with the latest reported values.
# Columns D:J hold the seven most recent annual periods (most recent first).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BCH")

# Row 8: Total Revenue
$ws.Range("D8").Value = 2765700
$ws.Range("E8").Value = 2810100
$ws.Range("F8").Value = 2792000
$ws.Range("G8").Value = 3007000
$ws.Range("H8").Value = 2595900
$ws.Range("I8").Value = 2459000
$ws.Range("J8").Value = 2207500

# Row 15: Others
$ws.Range("D15").Value = -51800
$ws.Range("E15").Value = -48900
$ws.Range("F15").Value = -43400
$ws.Range("G15").Value = -48200
$ws.Range("H15").Value = -40700
$ws.Range("I15").Value = -103300
$ws.Range("J15").Value = -103300

# Row 17: Total Operating Expenses
$ws.Range("D17").Value = 1303900
$ws.Range("E17").Value = 1470000
$ws.Range("F17").Value = 1445300
$ws.Range("G17").Value = 1544000
$ws.Range("H17").Value = 1361300
$ws.Range("I17").Value = 1286300
$ws.Range("J17").Value = 1133600

# Row 18: Operating Income or Loss
$ws.Range("D18").Value = 1461900
$ws.Range("E18").Value = 1340100
$ws.Range("F18").Value = 1346600
$ws.Range("G18").Value = 1463000
$ws.Range("H18").Value = 1234700
$ws.Range("I18").Value = 1172600
$ws.Range("J18").Value = 1073900

# Row 20: Total Other Income/Expenses Net
$ws.Range("D20").Value = -446000
$ws.Range("E20").Value = -397400
$ws.Range("F20").Value = -434200
$ws.Range("G20").Value = -470500
$ws.Range("H20").Value = -294900
$ws.Range("I20").Value = -373200
$ws.Range("J20").Value = -333700

# Row 21: Earnings Before Interest And Taxes
$ws.Range("D21").Value = 1067700
$ws.Range("E21").Value = 991600
$ws.Range("F21").Value = 955900
$ws.Range("G21").Value = 1040700
$ws.Range("H21").Value = 980400
$ws.Range("I21").Value = 851100
$ws.Range("J21").Value = "NA"

# Row 23: Income Before Tax
$ws.Range("D23").Value = 1015800
$ws.Range("E23").Value = 942700
$ws.Range("F23").Value = 912500
$ws.Range("G23").Value = 992500
$ws.Range("H23").Value = 939700
$ws.Range("I23").Value = 799400
$ws.Range("J23").Value = 740300

# Row 24: Income Tax Expense
$ws.Range("D24").Value = 169100
$ws.Range("E24").Value = 130900
$ws.Range("F24").Value = 90700
$ws.Range("G24").Value = 117100
$ws.Range("H24").Value = 131000
$ws.Range("I24").Value = 94000
$ws.Range("J24").Value = 96200

# Row 26: Income After Tax
$ws.Range("D26").Value = 846700
$ws.Range("E26").Value = 811800
$ws.Range("F26").Value = 821700
$ws.Range("G26").Value = 875400
$ws.Range("H26").Value = 808800
$ws.Range("I26").Value = 705400
$ws.Range("J26").Value = 644100

# Row 27: Net Income From Continuing Ops
$ws.Range("D27").Value = 846700
$ws.Range("E27").Value = 811800
$ws.Range("F27").Value = 821700
$ws.Range("G27").Value = 875400
$ws.Range("H27").Value = 808800
$ws.Range("I27").Value = 705400
$ws.Range("J27").Value = 644100

# Row 32: Other Items
$ws.Range("D32").Value = 446000
$ws.Range("E32").Value = 397400
$ws.Range("F32").Value = 434200
$ws.Range("G32").Value = 470500
$ws.Range("H32").Value = 294900
$ws.Range("I32").Value = 373200
$ws.Range("J32").Value = 333700

# Row 33: Net Income
$ws.Range("D33").Value = 846700
$ws.Range("E33").Value = 811800
$ws.Range("F33").Value = 821700
$ws.Range("G33").Value = 875400
$ws.Range("H33").Value = 808800
$ws.Range("I33").Value = 705400
$ws.Range("J33").Value = 644100

# Row 35: Net Income Applicable To Common Shares
$ws.Range("D35").Value = 846700
$ws.Range("E35").Value = 811800
$ws.Range("F35").Value = 821700
$ws.Range("G35").Value = 875400
$ws.Range("H35").Value = 808800
$ws.Range("I35").Value = 705400
$ws.Range("J35").Value = 644100

# Row 41: Cash And Cash Equivalents
$ws.Range("D41").Value = 2891200
$ws.Range("E41").Value = 3429700
$ws.Range("F41").Value = 4825200
$ws.Range("G41").Value = 3631800
$ws.Range("H41").Value = 3286000
$ws.Range("I41").Value = 1462700
$ws.Range("J41").Value = 1844500

# Row 42: Short Term Investments
$ws.Range("D42").Value = 4892500
$ws.Range("E42").Value = 3066000
$ws.Range("F42").Value = 1341800
$ws.Range("G42").Value = 2070300
$ws.Range("H42").Value = 1152500
$ws.Range("I42").Value = 2261000
$ws.Range("J42").Value = 1420400

# Row 47: Long Term Investments
$ws.Range("D47").Value = 55900
$ws.Range("E47").Value = 47900
$ws.Range("F47").Value = 41300
$ws.Range("G47").Value = 37200
$ws.Range("H47").Value = 21200
$ws.Range("I47").Value = 17200
$ws.Range("J47").Value = 19400

# Row 48: Property Plant and Equipment
$ws.Range("D48").Value = 338900
$ws.Range("E48").Value = 343600
$ws.Range("F48").Value = 339100
$ws.Range("G48").Value = 1033100
$ws.Range("H48").Value = 314400
$ws.Range("I48").Value = 326200
$ws.Range("J48").Value = 330700

# Row 49: Goodwill
$ws.Range("D49").Value = 57400
$ws.Range("E49").Value = 43100
$ws.Range("F49").Value = 39300
$ws.Range("G49").Value = 73700
$ws.Range("H49").Value = 106200
$ws.Range("I49").Value = 111100
$ws.Range("J49").Value = 119100

# Row 52: Other Assets
$ws.Range("D52").Value = 393100
$ws.Range("E52").Value = 414100
$ws.Range("F52").Value = 376300
$ws.Range("G52").Value = 298200
$ws.Range("H52").Value = 82900
$ws.Range("I52").Value = 82000
$ws.Range("J52").Value = 88200

# Row 54: Total Assets
$ws.Range("D54").Value = 48251600
$ws.Range("E54").Value = 46354500
$ws.Range("F54").Value = 46000600
$ws.Range("G54").Value = 40639400
$ws.Range("H54").Value = 38034400
$ws.Range("I54").Value = 34101400
$ws.Range("J54").Value = 31981300

# Row 57: Accounts Payable
$ws.Range("D57").Value = 714200
$ws.Range("E57").Value = 501900
$ws.Range("F57").Value = 534000
$ws.Range("G57").Value = 320900
$ws.Range("H57").Value = 223400
$ws.Range("I57").Value = 106800
$ws.Range("J57").Value = 228500

# Row 58: Short/Current Long Term Debt
$ws.Range("D58").Value = 3693500
$ws.Range("G58").Value = 1687500
$ws.Range("H58").Value = 1438600

# Row 59: Other Current Liabilities
$ws.Range("D59").Value = 466800
$ws.Range("E59").Value = 421400
$ws.Range("F59").Value = 43000
$ws.Range("G59").Value = 154400
$ws.Range("H59").Value = 110400
$ws.Range("I59").Value = 34100
$ws.Range("J59").Value = 4500

# Row 61: Long Term Debt
$ws.Range("D61").Value = 8090800
$ws.Range("E61").Value = 10884100
$ws.Range("F61").Value = 11473200
$ws.Range("G61").Value = 9324600
$ws.Range("H61").Value = 7122800
$ws.Range("I61").Value = 7239800
$ws.Range("J61").Value = 6900100

# Row 62: Other Liabilities
$ws.Range("D62").Value = 563000
$ws.Range("E62").Value = 553900
$ws.Range("F62").Value = 987800
$ws.Range("G62").Value = 936000
$ws.Range("H62").Value = 227300
$ws.Range("I62").Value = 94900
$ws.Range("J62").Value = 89100

# Row 66: Total Liabilities
$ws.Range("D66").Value = 43686200
$ws.Range("E66").Value = 42110000
$ws.Range("F66").Value = 41972700
$ws.Range("G66").Value = 36912700
$ws.Range("H66").Value = 34096200
$ws.Range("I66").Value = 30638800
$ws.Range("J66").Value = 28981500

# Row 72: Retained Earnings
$ws.Range("D72").Value = 1242200
$ws.Range("E72").Value = 1140800
$ws.Range("F72").Value = 942600
$ws.Range("G72").Value = 802800
$ws.Range("H72").Value = 1185800
$ws.Range("I72").Value = 1235300
$ws.Range("J72").Value = 1077500

# Row 76: Total Stockholder Equity
$ws.Range("D76").Value = 4565400
$ws.Range("E76").Value = 4244500
$ws.Range("F76").Value = 4027900
$ws.Range("G76").Value = 3726700
$ws.Range("H76").Value = 3938200
$ws.Range("I76").Value = 3462500
$ws.Range("J76").Value = 2999800

# Row 81: Net Income
$ws.Range("D81").Value = 846700
$ws.Range("E81").Value = 811800
$ws.Range("F81").Value = 821700
$ws.Range("G81").Value = 875400
$ws.Range("H81").Value = 808800
$ws.Range("I81").Value = 705400
$ws.Range("J81").Value = 644100

# Row 83: Depreciation
$ws.Range("D83").Value = 51800
$ws.Range("E83").Value = 48900
$ws.Range("F83").Value = 43400
$ws.Range("G83").Value = 48200
$ws.Range("H83").Value = 40700
$ws.Range("I83").Value = 51700
$ws.Range("J83").Value = "NA"

# Row 89: Total Cash Flow From Operating Activities
$ws.Range("D89").Value = 1569200
$ws.Range("E89").Value = 500400
$ws.Range("F89").Value = -2092900
$ws.Range("G89").Value = -337400
$ws.Range("H89").Value = -205500
$ws.Range("I89").Value = -611500
$ws.Range("J89").Value = -668700

# Row 91: Capital Expenditures
$ws.Range("D91").Value = -34100
$ws.Range("E91").Value = -40900
$ws.Range("F91").Value = -46300
$ws.Range("G91").Value = -46300
$ws.Range("H91").Value = -18000
$ws.Range("I91").Value = -26400
$ws.Range("J91").Value = -32400

# Row 94: Total Cash Flows From Investing Activities
$ws.Range("D94").Value = -1733600
$ws.Range("E94").Value = 592700
$ws.Range("F94").Value = 588100
$ws.Range("G94").Value = 72400
$ws.Range("H94").Value = -601600
$ws.Range("I94").Value = 282800
$ws.Range("J94").Value = "NA"

# Row 96: Dividends Paid
$ws.Range("D96").Value = -502800
$ws.Range("E96").Value = -539000
$ws.Range("F96").Value = -540100
$ws.Range("G96").Value = -541100
$ws.Range("H96").Value = -504900
$ws.Range("I96").Value = -436300
$ws.Range("J96").Value = -410400

# Row 100: Total Cash Flows From Financing Activities
$ws.Range("D100").Value = 195000
$ws.Range("E100").Value = -1046000
$ws.Range("F100").Value = 1784300
$ws.Range("G100").Value = 555000
$ws.Range("H100").Value = 1133200
$ws.Range("I100").Value = 70700
$ws.Range("J100").Value = "NA"

# Row 101: Effect Of Exchange Rate Changes
$ws.Range("D101").Value = -56400
$ws.Range("E101").Value = -42500
$ws.Range("F101").Value = 114900
$ws.Range("G101").Value = 67900
$ws.Range("H101").Value = 88800
$ws.Range("I101").Value = -46600
$ws.Range("J101").Value = "NA"

# Row 102: Change In Cash and Cash Equivalents
$ws.Range("D102").Value = -25800
$ws.Range("E102").Value = 4500
$ws.Range("F102").Value = 394400
$ws.Range("G102").Value = 357900
$ws.Range("H102").Value = 414900
$ws.Range("I102").Value = -304700
$ws.Range("J102").Value = 190700
